# Apply the edit described by the diff:
#  - A10: number 21001 -> text "NA" (picks up column A's right-aligned style,
#         matching the style already used by the other cells in column A)
#  - E10: "removeNegative" -> new text "recodeBMIcon" (picks up the same
#         Menlo/tan style already used by E8 for "recodeHouse")
#  - Selection moves from B10 to A10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing formatting from sibling cells that already carry the
# target styles, so we reuse the workbook's existing style records instead
# of fabricating new ones.
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("E8").Copy()
$ws.Range("E10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Now set the new cell contents.
$ws.Range("A10").Value = "NA"
$ws.Range("E10").Value = "recodeBMIcon"

# Move the active selection to A10 (was B10).
$ws.Range("A10").Select()
